$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 88, shifting existing rows 88-134 down to 89-135
$ws.Rows.Item(88).Insert()

# Populate the new row 88 with the new record
$ws.Cells.Item(88, 1).Value = 9
$ws.Cells.Item(88, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(88, 3).Value = "Metropolitana"
$ws.Cells.Item(88, 4).Value = 44813
$ws.Cells.Item(88, 5).Value = 13
$ws.Cells.Item(88, 6).Value = 100112022
$ws.Cells.Item(88, 7).Value = "Arveja Verde"
$ws.Cells.Item(88, 8).Value = "Perfection"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 20
$ws.Cells.Item(88, 11).Value = 40000
$ws.Cells.Item(88, 12).Value = 40000
$ws.Cells.Item(88, 13).Value = 40000
$ws.Cells.Item(88, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(88, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(88, 16).Value = 1600
$ws.Cells.Item(88, 17).Value = 25
$ws.Cells.Item(88, 18).Value = "Hortaliza"
